$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number, date range) ---
$ws.Range("A8").Value = "Volume 32   Number  37"
$ws.Range("C9").Value = "Report Covering the Week  9/8/2025  Through  9/14/2025"

# --- Crime complaints table updates (rows 14-33) ---
$ws.Range("C14").Value = 1
$ws.Range("C14").NumberFormat = "#,##0"
$ws.Range("D14").Value = 2
$ws.Range("D14").NumberFormat = "#,##0"
$ws.Range("E14").Value = -50
$ws.Range("E14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G14").Value = 2
$ws.Range("G14").NumberFormat = "#,##0"
$ws.Range("H14").Value = -50
$ws.Range("H14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I14").Value = 6
$ws.Range("J14").Value = 4
$ws.Range("K14").Value = 50
$ws.Range("L14").Value = 200
$ws.Range("M14").Value = 500
$ws.Range("N14").Value = -57.142857142857
$ws.Range("D15").Value = "0"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "***.*"
$ws.Range("E15").NumberFormat = "General"
$ws.Range("G15").Value = 2
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = 33.333333333333
$ws.Range("F16").Value = 30
$ws.Range("G16").Value = 33
$ws.Range("H16").Value = -9.090909090909
$ws.Range("I16").Value = 225
$ws.Range("J16").Value = 312
$ws.Range("K16").Value = -27.884615384615
$ws.Range("L16").Value = -2.173913043478
$ws.Range("M16").Value = -8.163265306122
$ws.Range("N16").Value = -76.240760295670
$ws.Range("C17").Value = 14
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = 7.692307692307
$ws.Range("F17").Value = 37
$ws.Range("G17").Value = 48
$ws.Range("H17").Value = -22.916666666666
$ws.Range("I17").Value = 389
$ws.Range("J17").Value = 470
$ws.Range("K17").Value = -17.234042553191
$ws.Range("L17").Value = 12.427745664739
$ws.Range("M17").Value = 62.083333333333
$ws.Range("N17").Value = 23.492063492063
$ws.Range("C18").Value = 4
$ws.Range("E18").Value = -42.857142857142
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = -28.571428571428
$ws.Range("I18").Value = 114
$ws.Range("J18").Value = 165
$ws.Range("K18").Value = -30.909090909090
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -44.390243902439
$ws.Range("N18").Value = -91.988756148981
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = -33.333333333333
$ws.Range("F19").Value = 54
$ws.Range("G19").Value = 79
$ws.Range("H19").Value = -31.645569620253
$ws.Range("I19").Value = 515
$ws.Range("J19").Value = 687
$ws.Range("K19").Value = -25.036390101892
$ws.Range("L19").Value = -12.116040955631
$ws.Range("M19").Value = 51.470588235294
$ws.Range("N19").Value = -50.858778625954
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 400
$ws.Range("I20").Value = 175
$ws.Range("J20").Value = 208
$ws.Range("K20").Value = -15.865384615384
$ws.Range("L20").Value = -24.892703862660
$ws.Range("M20").Value = 8.024691358024
$ws.Range("N20").Value = -88.796414852752
$ws.Range("C21").Value = 44
$ws.Range("E21").Value = -6.382978723404
$ws.Range("F21").Value = 162
$ws.Range("G21").Value = 202
$ws.Range("H21").Value = -19.801980198019
$ws.Range("I21").Value = 1445
$ws.Range("J21").Value = 1872
$ws.Range("K21").Value = -22.809829059829
$ws.Range("L21").Value = -5.493786788750
$ws.Range("M21").Value = 19.126133553174
$ws.Range("N21").Value = -72.894391296192
$ws.Range("C22").Value = "0"
$ws.Range("C22").NumberFormat = "General"
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -100
$ws.Range("J22").Value = 62
$ws.Range("K22").Value = -43.548387096774
$ws.Range("L22").Value = -46.153846153846
$ws.Range("C24").Value = 33
$ws.Range("D24").Value = 38
$ws.Range("E24").Value = -13.157894736842
$ws.Range("F24").Value = 146
$ws.Range("G24").Value = 140
$ws.Range("H24").Value = 4.285714285714
$ws.Range("I24").Value = 1022
$ws.Range("J24").Value = 1507
$ws.Range("K24").Value = -32.183145321831
$ws.Range("L24").Value = -26.315789473684
$ws.Range("M24").Value = 30.357142857142
$ws.Range("C25").Value = 8
$ws.Range("E25").Value = -42.857142857142
$ws.Range("F25").Value = 29
$ws.Range("G25").Value = 59
$ws.Range("H25").Value = -50.847457627118
$ws.Range("I25").Value = 344
$ws.Range("J25").Value = 803
$ws.Range("K25").Value = -57.160647571606
$ws.Range("L25").Value = -51.066856330014
$ws.Range("C26").Value = 23
$ws.Range("D26").Value = 22
$ws.Range("E26").Value = 4.545454545454
$ws.Range("F26").Value = 83
$ws.Range("G26").Value = 114
$ws.Range("H26").Value = -27.192982456140
$ws.Range("I26").Value = 758
$ws.Range("J26").Value = 977
$ws.Range("K26").Value = -22.415557830092
$ws.Range("L26").Value = 10.818713450292
$ws.Range("M26").Value = 12.965722801788
$ws.Range("D27").Value = "0"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "***.*"
$ws.Range("E27").NumberFormat = "General"
$ws.Range("G27").Value = 3
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = -33.333333333333
$ws.Range("F28").Value = 15
$ws.Range("G28").Value = 13
$ws.Range("H28").Value = 15.384615384615
$ws.Range("I28").Value = 85
$ws.Range("J28").Value = 102
$ws.Range("K28").Value = -16.666666666666
$ws.Range("L28").Value = -26.086956521739
$ws.Range("C29").Value = 1
$ws.Range("C29").NumberFormat = "#,##0"
$ws.Range("D29").Value = 2
$ws.Range("D29").NumberFormat = "#,##0"
$ws.Range("E29").Value = -50
$ws.Range("E29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F29").Value = 1
$ws.Range("F29").NumberFormat = "#,##0"
$ws.Range("G29").Value = 2
$ws.Range("G29").NumberFormat = "#,##0"
$ws.Range("H29").Value = -50
$ws.Range("H29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I29").Value = 1
$ws.Range("I29").NumberFormat = "#,##0"
$ws.Range("J29").Value = 6
$ws.Range("K29").Value = -83.333333333333
$ws.Range("L29").Value = -80
$ws.Range("M29").Value = -83.333333333333
$ws.Range("N29").Value = -97.872340425531
$ws.Range("C30").Value = 1
$ws.Range("C30").NumberFormat = "#,##0"
$ws.Range("D30").Value = 2
$ws.Range("D30").NumberFormat = "#,##0"
$ws.Range("E30").Value = -50
$ws.Range("E30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F30").Value = 1
$ws.Range("F30").NumberFormat = "#,##0"
$ws.Range("G30").Value = 2
$ws.Range("G30").NumberFormat = "#,##0"
$ws.Range("H30").Value = -50
$ws.Range("H30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I30").Value = 1
$ws.Range("I30").NumberFormat = "#,##0"
$ws.Range("J30").Value = 4
$ws.Range("K30").Value = -75
$ws.Range("L30").Value = -75
$ws.Range("M30").Value = -80
$ws.Range("N30").Value = -97.674418604651
$ws.Range("L33").Value = -80
